$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update virus_full_name for existing rows 3-6 to use sequential numbering
$ws.Range("D3").Value = " Carnivore amdoparvovirus 2"
$ws.Range("D4").Value = " Carnivore amdoparvovirus 3"
$ws.Range("D5").Value = " Carnivore amdoparvovirus 4"
$ws.Range("D6").Value = " Carnivore amdoparvovirus 5"

# Update row 7 (Labrador amdoparvovirus) to new naming scheme
$ws.Range("B7").Value = "LaAV-1"
$ws.Range("D7").Value = "Labrador-amdoparvovirus-1"
$ws.Range("E7").Value = "NA"
$ws.Range("O7").Value = "NULL"

# Add new row 8 - Rattus nitidus parvovirus
$ws.Range("A8").Value = "KJ641663"
$ws.Range("B8").Value = "RtRn-ParV"
$ws.Range("C8").Value = "Parvovirinae"
$ws.Range("D8").Value = "Rattus nitidus parvovirus"
$ws.Range("E8").Value = "NA"
$ws.Range("F8").Value = "Rattus nitidus"
$ws.Range("G8").Value = "Amdoparvovirus"
$ws.Range("H8").Value = "NK"
$ws.Range("I8").Value = "NK"
$ws.Range("J8").Value = "NK"
$ws.Range("K8").Value = "NK"
$ws.Range("L8").Value = "NK"
$ws.Range("M8").Value = "NK"
$ws.Range("N8").Value = "NK"
$ws.Range("O8").Value = "NULL"

# Add new row 9 - Rhinolophus lepidus parvovirus
$ws.Range("A9").Value = "KY432922"
$ws.Range("B9").Value = "BtRl-PV"
$ws.Range("C9").Value = "Parvovirinae"
$ws.Range("D9").Value = "Rhinolophus lepidus parvovirus"
$ws.Range("E9").Value = "NA"
$ws.Range("F9").Value = "Rhinolophus lepidus"
$ws.Range("G9").Value = "Amdoparvovirus"
$ws.Range("H9").Value = "NK"
$ws.Range("I9").Value = "NK"
$ws.Range("J9").Value = "NK"
$ws.Range("K9").Value = "NK"
$ws.Range("L9").Value = "NK"
$ws.Range("M9").Value = "NK"
$ws.Range("N9").Value = "NK"
$ws.Range("O9").Value = "NULL"

# Update selection to match target state
$ws.Range("B9").Select()
